# Commit: "test push avec --set-upstream"
#
# The two list-paragraphs
#   "J'ajoute du texte"
#   "Je test pour voir si changement sur github"
# become a single list-paragraph:
#   "Test avec –set-upstream."
# (the trailing <w:bookmarkStart/bookmarkEnd name="_GoBack"> that used to sit
# at the end of the second paragraph now sits at the end of the merged one).

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $text) {
    foreach ($p in $doc.Paragraphs) {
        $t = $p.Range.Text
        if ($t.Length -gt 0) {
            $lastCode = [int][char]$t.Substring($t.Length - 1, 1)
            if ($lastCode -eq 13 -or $lastCode -eq 7) {
                $t = $t.Substring(0, $t.Length - 1)
            }
        }
        if ($t -eq $text) {
            return $p
        }
    }
    return $null
}

$oldText1 = "J" + [char]0x2019 + "ajoute du texte"
$oldText2 = "Je test pour voir si changement sur github"
$newText  = "Test avec " + [char]0x2013 + "set-upstream."

# 1) Swap the text of the first paragraph for the new sentence (keeps the
#    paragraph's own pPr/run rPr - same run, same formatting).
$p1 = Find-ParagraphByText $d $oldText1
if ($p1 -eq $null) {
    throw "Could not find paragraph containing '$oldText1'"
}
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Text = $newText

# 2) Remove the now-obsolete text of the second paragraph, leaving its
#    paragraph mark and the _GoBack bookmark that follows it untouched.
$p2 = Find-ParagraphByText $d $oldText2
if ($p2 -eq $null) {
    throw "Could not find paragraph containing '$oldText2'"
}
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$r2.Delete()

# 3) Merge the (now empty) second paragraph back into the first one by
#    deleting the first paragraph's trailing mark, so the bookmark ends up
#    right after the new run, inside the merged paragraph.
$p1b = Find-ParagraphByText $d $newText
if ($p1b -eq $null) {
    throw "Could not relocate the edited paragraph"
}
$mark = $d.Range($p1b.Range.End - 1, $p1b.Range.End)
$mark.Delete()

Write-Output "Merged the two list paragraphs into: $newText"
